# chore: simulator full-month coverage, persist logs, fix employees
#
# Updates the Boban Abbate 2026-01-19 weekly export:
#  - fixes client names on both the "Weekly Timesheet" and "Jason Schema"
#    sheets (simulator now emits a full month of clients, the old names
#    were placeholders)
#  - fills in the simulator-computed Rate / Total figures that were
#    previously zeroed out
#  - fixes the employee id that was logged for this run

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")
$ws2 = $wb.Worksheets.Item("Jason Schema")

# --- client name fixes -----------------------------------------------
# row 2: Evans   -> Schauer
# row 3: Oglesby -> Muncey
# row 4: Muncey  -> Moulton
# row 5: Lucas   -> Regan
# row 6: Bailey  -> Hendricks
$newNames = @("Schauer", "Muncey", "Moulton", "Regan", "Hendricks")

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $newNames[$i]
    $ws2.Cells.Item($row, 4).Value = $newNames[$i]
}

# --- employee id fix ---------------------------------------------------
for ($row = 2; $row -le 6; $row++) {
    $ws2.Cells.Item($row, 2).Value = "emp_ga4rqytu"
}

# --- Weekly Timesheet: rate / total for each daily row -----------------
for ($row = 2; $row -le 6; $row++) {
    $ws1.Cells.Item($row, 5).Value = 110
    $ws1.Cells.Item($row, 6).Value = 880
}

# --- Weekly Timesheet: rolled-up totals ---------------------------------
$ws1.Range("F8").Value = 4400   # SUBTOTAL
$ws1.Range("F11").Value = 4400  # HOURLY SUBTOTAL
$ws1.Range("F13").Value = 4400  # GRAND TOTAL

# --- Jason Schema: rate / total for each daily row ----------------------
for ($row = 2; $row -le 6; $row++) {
    $ws2.Cells.Item($row, 6).Value = 110
    $ws2.Cells.Item($row, 7).Value = 880
}
